$d = $word.ActiveDocument

$replacements = @(
    @("43×27=", "26×78="),
    @("56×49=", "59×92="),
    @("66×55=", "73×43="),
    @("20×81=", "48×22="),
    @("37×14=", "18×82="),
    @("65×87=", "28×82="),
    @("69×18=", "47×51="),
    @("75×64=", "17×78="),
    @("95×97=", "15×11="),
    @("53×42=", "48×68="),
    @("95×60=", "49×73="),
    @("12×85=", "51×35="),
    @("58×23=", "71×90="),
    @("86×97=", "16×40="),
    @("14×37=", "34×48="),
    @("16×67=", "52×57="),
    @("92×52=", "74×71="),
    @("59×40=", "91×51="),
    @("20×49=", "49×69="),
    @("22×54=", "18×98="),
    @("22×58=", "80×12="),
    @("64×19=", "84×89="),
    @("26×24=", "69×95="),
    @("26×52=", "79×11="),
    @("46×70=", "23×66=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
